$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells remain plain text (matches the original inlineStr cell type)
# by forcing a text number-format before assigning string values that could
# otherwise be auto-converted to numbers by Excel's input parser.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.286.70"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.431.80"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.00"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.57"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.17%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.33%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.432.82"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.86%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.08%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.15%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.28"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.64%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.51"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.30%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.40%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.154.18"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.431.22"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.26"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "323.91"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.86"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.94%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.80%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.44"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.62%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.73"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.76"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.69%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "553.18"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.49%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.53%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0951"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.40%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.42"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.37%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.17%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.08%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.84%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.09%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.77%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "150.31"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.67%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.80%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "147.74"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.80%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0533"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.27"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.598"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.60%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.88%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.99%  "

# Rows 38/39: PolygonEcosystemToken and RenderToken swap list positions,
# each also carrying freshly updated price/volume figures.
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.382"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.22%  "

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.56"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.07%  "
